# Updates odds figures in the "Jogos do Dia" Betfair Back/Lay sheet to
# reflect the latest market prices (Atualizando o arquivo XLSX).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("J2").Value = 3.8
$ws.Range("K2").Value = 3.85
$ws.Range("P2").Value = 1.81
$ws.Range("AI2").Value = 85
$ws.Range("AL2").Value = 42

# Row 3
$ws.Range("I3").Value = 1.81
$ws.Range("K3").Value = 3.75
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.79
$ws.Range("Q3").Value = 2.24
$ws.Range("U3").Value = 1.88
$ws.Range("AK3").Value = 90
$ws.Range("AM3").Value = 160

# Row 4
$ws.Range("K4").Value = 3.6
$ws.Range("O4").Value = 1.48

# Row 5
$ws.Range("G5").Value = 5.4
$ws.Range("O5").Value = 1.35
$ws.Range("P5").Value = 1.91
$ws.Range("U5").Value = 1.97
$ws.Range("AB5").Value = 17.5

# Row 6
$ws.Range("F6").Value = 2.86
$ws.Range("H6").Value = 2.62
$ws.Range("I6").Value = 2.72
$ws.Range("L6").Value = 1.39
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 1.84
$ws.Range("S6").Value = 3.65
$ws.Range("T6").Value = 1.76
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.58
$ws.Range("W6").Value = 1.52
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 11
$ws.Range("AF6").Value = 1000
$ws.Range("AH6").Value = 17.5
$ws.Range("AI6").Value = 44
$ws.Range("AN6").Value = 30
$ws.Range("AO6").Value = 26

# Row 7
$ws.Range("P7").Value = 1.97
$ws.Range("Q7").Value = 2
$ws.Range("Z7").Value = 15
$ws.Range("AO7").Value = 19.5

# Row 8
$ws.Range("F8").Value = 1.64
$ws.Range("G8").Value = 1.65
$ws.Range("H8").Value = 6.2
$ws.Range("I8").Value = 6.4
$ws.Range("J8").Value = 4.3
$ws.Range("K8").Value = 4.4
$ws.Range("V8").Value = 1.18
$ws.Range("W8").Value = 2.52
$ws.Range("AB8").Value = 11
$ws.Range("AG8").Value = 9.800000000000001
$ws.Range("AL8").Value = 36
